# "Layout atualizado separando tabela de resultados dos graficos."
#
# 1) Remove the (now redundant) "Parametros_Modular" sheet - its content was
#    a duplicate/variant of "Parametros" used for a modular-cost experiment
#    that is no longer part of the layout.
# 2) Re-point the view state: Configs becomes the active tab, and a handful
#    of sheets get their remembered selection moved to a different cell.
# 3) On "Custos", drop the helper formulas in D7:D11 (Iniciativa2) and bump
#    the Iniciativa1 / Iniciativa2 / TodasIniciativas cost figures up by two
#    orders of magnitude (they were placeholder numbers), widening column A
#    to fit the longer labels now that the sheet is presented on its own.

$wb = $excel.ActiveWorkbook

# --- 1. Delete the obsolete "Parametros_Modular" sheet -----------------
$wsModular = $wb.Worksheets.Item("Parametros_Modular")
$wsModular.Delete()

# --- 2. Update the Custos ("Custos") figures ----------------------------
$wsCustos = $wb.Worksheets.Item("Custos")

# Iniciativa1 (rows 2-6): 500 -> 50000
$wsCustos.Range("D2:D6").Value = 50000

# Iniciativa2 (rows 7-11): was "=D2*2" (shared formula) -> flat 800000
$wsCustos.Range("D7:D11").Value = 800000

# TodasIniciativas (rows 17-21): 1500 -> 150000
$wsCustos.Range("D17:D21").Value = 150000

# Column A is now wider to fit its labels without the bestFit flag
$wsCustos.Columns.Item(1).ColumnWidth = 18.15

# --- 3. Update remembered selections on each sheet ----------------------
$wsLista = $wb.Worksheets.Item("Lista_de_Parâmetros")
$wsLista.Activate()
$wsLista.Range("C9").Select()

$wsDados = $wb.Worksheets.Item("Dados_Projetados")
$wsDados.Activate()
$wsDados.Range("A4").Select()

$wsParametros = $wb.Worksheets.Item("Parametros")
$wsParametros.Activate()
$wsParametros.Range("A1").Select()

$wsCustos.Activate()
$wsCustos.Range("D7").Select()

# --- 4. Make "Configs" the active tab (tabSelected/activeTab) ----------
$wsConfigs = $wb.Worksheets.Item("Configs")
$wsConfigs.Activate()
$wsConfigs.Range("A1").Select()
